$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new data row (row 4) with double-byte character test data
# (shared-string insertion order matters: D4, then A4, then B4)
$ws.Range("D4").Value = "今日は晴れです。"
$ws.Range("A4").Value = "DoubleByteCharacters.001"
$ws.Range("B4").Value = "今日の天気を教えてください."

# C4 stays empty but keeps the wrap-text style used throughout column D
$ws.Range("C4").WrapText = $true
$ws.Range("D4").WrapText = $true

# Match row 3's row height (17pt) for the new row
$ws.Rows.Item(4).RowHeight = 17

$ws.Range("F8").Select()
